$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Employees"

$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"

$ws.Range("A2").Value = "Mary"
$ws.Range("A3").Value = "Vinod"
$ws.Range("A4").Value = "Mansoor"
$ws.Range("A5").Value = "Linda"

$ws.Range("B2").Value = "Brown"
$ws.Range("B3").Value = "Kumar"
$ws.Range("B4").Value = "Khan"
$ws.Range("B5").Value = "Lee"

$ws.Range("C1").Value = "Job_Title"
$ws.Range("D1").Value = "ID"

$ws.Range("C2").Value = "PO"
$ws.Range("C3").Value = "BA"
$ws.Range("C4").Value = "QA"
$ws.Range("C5").Value = "Developer"

$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4

$ws.Range("C4").Select()
